# Generate Report for Handback
#
# The handback for e0d4b6a4-8890-41c2-b163-4f3c1f04c079 came back on a stale
# source revision, so the zh-cn / de-de status rows for that file now show
# the target xlf name, the "not latest" error, and the handback timestamps
# that the CI run produced.

$wb = $excel.ActiveWorkbook

$targetMd   = "e0d4b6a4-8890-41c2-b163-4f3c1f04c079.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/252aec560098b1f9d704ff46a520338291625667/e2e/e0d4b6a4-8890-41c2-b163-4f3c1f04c079.md"
$staleMsg   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4633feffe92a5ddbbffff677c648adaaaa2135f/e2e/e0d4b6a4-8890-41c2-b163-4f3c1f04c079.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/252aec560098b1f9d704ff46a520338291625667/e2e/e0d4b6a4-8890-41c2-b163-4f3c1f04c079.md."

# --- zh-cn sheet, row 7 (e0d4b6a4-8890-41c2-b163-4f3c1f04c079) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value2 = "e0d4b6a4-8890-41c2-b163-4f3c1f04c079.1bc31c5a452f1ead8f88474d4c5649ecd0c4982f.zh-cn.xlf"
$wsZh.Range("K7").Value2 = $staleMsg
$wsZh.Range("P7").Value2 = "2016-09-01 21:03:34"

# Latest Target File (I7) becomes a link back to the source .md, same as
# every other populated row in the column.
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetMd) | Out-Null

# --- de-de sheet, row 7 (e0d4b6a4-8890-41c2-b163-4f3c1f04c079) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value2 = "e0d4b6a4-8890-41c2-b163-4f3c1f04c079.1bc31c5a452f1ead8f88474d4c5649ecd0c4982f.de-de.xlf"
$wsDe.Range("K7").Value2 = "2016-09-01 21:03:58"
$wsDe.Range("P7").Value2 = "2016-09-01 21:03:34"

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetMd) | Out-Null
